$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OperationScenario_Component_Reg")

# norm_outside_temperature for region DE, year 2019 changed from -12 to -10
$ws.Range("D2").Value = -10
